# Updates cryptos list values (price & 1h volume change) per the
# "Updated cryptos list" GitHub Actions commit.
# Also fixes the ordering of rows 15/16 (Polkadot <-> WrappedliquidstakedEther2.0)
# which had swapped in the source data.
#
# Numeric-looking price/percentage values are plain text in this sheet
# (stored as inline strings), so each assignment is prefixed with a leading
# apostrophe to force Excel to keep them as text instead of auto-converting
# them to numbers (which would lose formatting like "71.257.40" or drop
# trailing zeros such as "591.58" -> 591.58000000000004).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'71.257.40"
$ws.Range('E2').Value = "'  +2.66%  "
$ws.Range('D3').Value = "'3.524.90"
$ws.Range('E3').Value = "'  +3.94%  "
$ws.Range('D4').Value = "'0.998"
$ws.Range('E4').Value = "'  -0.17%  "
$ws.Range('D5').Value = "'591.58"
$ws.Range('E5').Value = "'  +0.60%  "
$ws.Range('D6').Value = "'180.30"
$ws.Range('E6').Value = "'  -0.17%  "
$ws.Range('D7').Value = "'3.499.65"
$ws.Range('E7').Value = "'  +3.42%  "
$ws.Range('D8').Value = "'0.605"
$ws.Range('E8').Value = "'  +1.40%  "
$ws.Range('E9').Value = "'  -0.04%  "
$ws.Range('D10').Value = "'0.208"
$ws.Range('E10').Value = "'  +6.53%  "
$ws.Range('D11').Value = "'0.599"
$ws.Range('E11').Value = "'  +1.35%  "
$ws.Range('D12').Value = "'49.67"
$ws.Range('E12').Value = "'  +2.33%  "
$ws.Range('D13').Value = "'0.0000288"
$ws.Range('E13').Value = "'  +2.10%  "
$ws.Range('D14').Value = "'696.93"
$ws.Range('E14').Value = "'  +2.84%  "
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = "'4.045.81"
$ws.Range('E15').Value = "'  +2.81%  "
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').Value = "'8.82"
$ws.Range('E16').Value = "'  +2.33%  "
$ws.Range('D17').Value = "'71.070.24"
$ws.Range('E17').Value = "'  +2.31%  "
$ws.Range('D18').Value = "'3.481.32"
$ws.Range('E18').Value = "'  +2.74%  "
$ws.Range('E19').Value = "'  +1.20%  "
$ws.Range('D20').Value = "'18.13"
$ws.Range('E20').Value = "'  +2.75%  "
$ws.Range('D21').Value = "'11.58"
$ws.Range('E21').Value = "'  +2.57%  "
$ws.Range('D22').Value = "'0.920"
$ws.Range('E22').Value = "'  +1.74%  "
$ws.Range('D23').Value = "'5.52"
$ws.Range('E23').Value = "'  +1.60%  "
$ws.Range('D24').Value = "'17.40"
$ws.Range('E24').Value = "'  +1.08%  "
$ws.Range('D25').Value = "'102.69"
$ws.Range('E25').Value = "'  -0.82%  "
$ws.Range('D26').Value = "'3.99"
$ws.Range('E26').Value = "'  +1.22%  "
$ws.Range('D27').Value = "'2.74"
$ws.Range('E27').Value = "'  +0.58%  "
$ws.Range('D28').Value = "'9.83"
$ws.Range('E28').Value = "'  +1.56%  "
$ws.Range('D29').Value = "'34.50"
$ws.Range('E29').Value = "'  +1.07%  "
$ws.Range('D30').Value = "'8.97"
$ws.Range('E30').Value = "'  +3.10%  "
$ws.Range('D31').Value = "'7.33"
$ws.Range('E31').Value = "'  +4.11%  "
$ws.Range('D32').Value = "'4.02"
$ws.Range('E32').Value = "'  +11.59%  "
$ws.Range('D33').Value = "'580.37"
$ws.Range('E33').Value = "'  +4.42%  "
$ws.Range('D34').Value = "'11.21"
$ws.Range('E34').Value = "'  +0.80%  "
$ws.Range('D35').Value = "'59.25"
$ws.Range('E35').Value = "'  +1.50%  "
$ws.Range('D36').Value = "'0.105"
$ws.Range('E36').Value = "'  -1.16%  "
$ws.Range('E37').Value = "'  -0.02%  "
$ws.Range('D38').Value = "'3.634.50"
$ws.Range('E38').Value = "'  -1.48%  "
$ws.Range('E39').Value = "'  +2.03%  "
$ws.Range('D40').Value = "'35.82"
$ws.Range('E40').Value = "'  +1.45%  "
$ws.Range('D41').Value = "'3.45"
$ws.Range('E41').Value = "'  +5.60%  "
$ws.Range('D42').Value = "'0.0₃0751"
$ws.Range('E42').Value = "'  +7.61%  "
$ws.Range('D43').Value = "'2.79"
$ws.Range('E43').Value = "'  +2.45%  "
$ws.Range('D44').Value = "'0.344"
$ws.Range('E44').Value = "'  +1.36%  "
$ws.Range('D45').Value = "'0.0433"
$ws.Range('D46').Value = "'3.36"
$ws.Range('E46').Value = "'  +2.45%  "
$ws.Range('D47').Value = "'2.75"
$ws.Range('E47').Value = "'  +2.83%  "
$ws.Range('D48').Value = "'1.47"
$ws.Range('E48').Value = "'  +4.47%  "
$ws.Range('D49').Value = "'0.131"
$ws.Range('E49').Value = "'  +0.85%  "
$ws.Range('E50').Value = "'  -0.25%  "
$ws.Range('D51').Value = "'134.32"
$ws.Range('E51').Value = "'  +0.78%  "

Write-Output "Applied cryptos list update."
